$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add two new task rows (17 and 18)
$ws.Range("A17").Value = "DB auto backups"
$ws.Range("B17").Value = "Ignas"
$ws.Range("C17").Value = "vidutinis"

$ws.Range("A18").Value = "Taisyklės"
$ws.Range("B18").Value = "Ignas"
$ws.Range("C18").Value = "vidutinis"

# Update the active selection on the sheet view
[void]$ws.Range("A23").Select()
